{"js": "// Update the worksheet date and every two-digit-by-two-digit\n// multiplication problem/answer pair with the new values from the\n// latest generated output.\nconst replacements = [\n  [\"2025-11-25 Tuesday\", \"2025-11-26 Wednesday\"],\n  [\"37\u00d734=1258\", \"99\u00d757=5643\"],\n  [\"92\u00d778=7176\", \"48\u00d792=4416\"],\n  [\"42\u00d744=1848\", \"40\u00d760=2400\"],\n  [\"19\u00d721=399\", \"15\u00d711=165\"],\n  [\"29\u00d788=2552\", \"43\u00d766=2838\"],\n  [\"14\u00d753=742\", \"43\u00d786=3698\"],\n  [\"13\u00d746=598\", \"70\u00d713=910\"],\n  [\"97\u00d731=3007\", \"86\u00d797=8342\"],\n  [\"50\u00d711=550\", \"99\u00d756=5544\"],\n  [\"33\u00d795=3135\", \"25\u00d786=2150\"],\n  [\"55\u00d781=4455\", \"40\u00d775=3000\"],\n  [\"85\u00d798=8330\", \"72\u00d752=3744\"],\n  [\"79\u00d762=4898\", \"53\u00d795=5035\"],\n  [\"21\u00d717=357\", \"39\u00d746=1794\"],\n  [\"95\u00d729=2755\", \"66\u00d725=1650\"],\n  [\"99\u00d759=5841\", \"72\u00d793=6696\"],\n  [\"19\u00d746=874\", \"27\u00d797=2619\"],\n  [\"27\u00d794=2538\", \"61\u00d762=3782\"],\n  [\"22\u00d777=1694\", \"97\u00d773=7081\"],\n  [\"40\u00d720=800\", \"58\u00d782=4756\"],\n  [\"97\u00d744=4268\", \"54\u00d723=1242\"],\n  [\"97\u00d742=4074\", \"52\u00d780=4160\"],\n  [\"85\u00d730=2550\", \"79\u00d793=7347\"],\n  [\"69\u00d785=5865\", \"28\u00d743=1204\"],\n  [\"76\u00d750=3800\", \"18\u00d750=900\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every two-digit-by-two-digit\n# multiplication problem/answer pair with the new values from the\n# latest generated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-11-25 Tuesday', '2025-11-26 Wednesday'),\n    @('37\u00d734=1258', '99\u00d757=5643'),\n    @('92\u00d778=7176', '48\u00d792=4416'),\n    @('42\u00d744=1848', '40\u00d760=2400'),\n    @('19\u00d721=399', '15\u00d711=165'),\n    @('29\u00d788=2552', '43\u00d766=2838'),\n    @('14\u00d753=742', '43\u00d786=3698'),\n    @('13\u00d746=598', '70\u00d713=910'),\n    @('97\u00d731=3007', '86\u00d797=8342'),\n    @('50\u00d711=550', '99\u00d756=5544'),\n    @('33\u00d795=3135', '25\u00d786=2150'),\n    @('55\u00d781=4455', '40\u00d775=3000'),\n    @('85\u00d798=8330', '72\u00d752=3744'),\n    @('79\u00d762=4898', '53\u00d795=5035'),\n    @('21\u00d717=357', '39\u00d746=1794'),\n    @('95\u00d729=2755', '66\u00d725=1650'),\n    @('99\u00d759=5841', '72\u00d793=6696'),\n    @('19\u00d746=874', '27\u00d797=2619'),\n    @('27\u00d794=2538', '61\u00d762=3782'),\n    @('22\u00d777=1694', '97\u00d773=7081'),\n    @('40\u00d720=800', '58\u00d782=4756'),\n    @('97\u00d744=4268', '54\u00d723=1242'),\n    @('97\u00d742=4074', '52\u00d780=4160'),\n    @('85\u00d730=2550', '79\u00d793=7347'),\n    @('69\u00d785=5865', '28\u00d743=1204'),\n    @('76\u00d750=3800', '18\u00d750=900')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
